$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2..64 in column D hold fraction values (0-1) that need to be
# rescaled to percentage numbers (0-100), i.e. value * 100.
for ($r = 2; $r -le 64; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $cell.Value = $val * 100
    }
}
